# Update "问题记录及解决" sheet:
#  - Row2 (F): clarify the "ignore unity Library path" note text, date bumped.
#  - Row3: date bumped (content unchanged).
#  - Row4: date bumped; solution (E) now notes the watermark issue that was
#    previously split into the "supplement" (F) column; F gets new content
#    describing the VideoPlayer regression on Android.
#  - Row5: was an empty placeholder row (only D5 = "否"); now a fully
#    populated row documenting a new rendering-order / sprite-skin issue.
#  - Active selection moves from F4 to F2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Dates bumped on rows 2-4 ---
$ws.Range("A2").Value2 = 44981
$ws.Range("A3").Value2 = 44982
$ws.Range("A4").Value2 = 44982

# --- Row 4: solution (E) absorbs the old watermark note; F gets a fresh note ---
$ws.Range("E4").Value2 = "使用Avpro插件播放视频，但有水印，正在解决"
$ws.Range("F4").Value2 = "使用回VideoPlayer后安卓机上又出现视频了"

# --- Row 5 (was a near-empty row with only D5 = "否") now fully populated ---
# Pick up A4's date style (centered, numFmt 14) for the new A5 date cell.
$ws.Range("A4").Copy($ws.Range("A5"))
$ws.Range("A5").Value2 = 44984
$ws.Range("B5").Value2 = "技术问题`n渲染顺序"
$ws.Range("C5").Value2 = "使用的sprite skin，人物朝左，理应左手在上，右手遮挡部分不可见，然而无论怎么调顺序双手都在身体下面"
$ws.Range("D5").Value2 = "是"
$ws.Range("E5").Value2 = "在sprite editor面板修改骨骼的depth属性，即可改变骨骼间的层级"
$ws.Range("F5").Value2 = ""

# --- Row 2: clarify the "ignore unity Library path" note text (edited last) ---
$ws.Range("F2").Value2 = "忽略文件路径为，以所在目录为根目录的，相对路径"

# --- Active selection moves from F4 to F2 ---
$ws.Range("F2").Select()
